$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 100 and row 101 (columns B:AC)
$ws.Cells.Item(100, 2).Value = 6867461
$ws.Cells.Item(100, 6).Value = 'FK Zeleziarne Podbrezova'
$ws.Cells.Item(100, 7).Value = 'Slovan Bratislava'
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 9).Value = 6
$ws.Cells.Item(100, 10).Value = 'A'
$ws.Cells.Item(100, 11).Value = 3.25
$ws.Cells.Item(100, 12).Value = 3.5
$ws.Cells.Item(100, 13).Value = 2
$ws.Cells.Item(100, 14).Value = 3.5
$ws.Cells.Item(100, 15).Value = 3.75
$ws.Cells.Item(100, 16).Value = 2
$ws.Cells.Item(100, 17).Value = 0.5
$ws.Cells.Item(100, 18).Value = 1.8
$ws.Cells.Item(100, 19).Value = 2
$ws.Cells.Item(100, 20).Value = 3
$ws.Cells.Item(100, 21).Value = 1.95
$ws.Cells.Item(100, 22).Value = 1.85
$ws.Cells.Item(100, 23).Value = -1
$ws.Cells.Item(100, 24).Value = -1
$ws.Cells.Item(100, 25).Value = 1
$ws.Cells.Item(100, 26).Value = -1
$ws.Cells.Item(100, 27).Value = 1
$ws.Cells.Item(100, 28).Value = 0.95
$ws.Cells.Item(100, 29).Value = -1
$ws.Cells.Item(101, 2).Value = 6867460
$ws.Cells.Item(101, 6).Value = 'FC Vion Zlate Moravce'
$ws.Cells.Item(101, 7).Value = 'FC Kosice'
$ws.Cells.Item(101, 8).Value = 1
$ws.Cells.Item(101, 9).Value = 1
$ws.Cells.Item(101, 10).Value = 'D'
$ws.Cells.Item(101, 11).Value = 2.3
$ws.Cells.Item(101, 12).Value = 3.3
$ws.Cells.Item(101, 13).Value = 2.875
$ws.Cells.Item(101, 14).Value = 2.75
$ws.Cells.Item(101, 15).Value = 3.1
$ws.Cells.Item(101, 16).Value = 2.75
$ws.Cells.Item(101, 17).Value = 0
$ws.Cells.Item(101, 18).Value = 1.875
$ws.Cells.Item(101, 19).Value = 1.925
$ws.Cells.Item(101, 20).Value = 2
$ws.Cells.Item(101, 21).Value = 1.85
$ws.Cells.Item(101, 22).Value = 1.95
$ws.Cells.Item(101, 23).Value = -1
$ws.Cells.Item(101, 24).Value = 2.1
$ws.Cells.Item(101, 25).Value = -1
$ws.Cells.Item(101, 26).Value = 0
$ws.Cells.Item(101, 27).Value = 0
$ws.Cells.Item(101, 28).Value = 0
$ws.Cells.Item(101, 29).Value = 0

# Swap row 124 and row 125 (columns B:AC)
$ws.Cells.Item(124, 2).Value = 6867489
$ws.Cells.Item(124, 6).Value = 'FC Spartak Trnava'
$ws.Cells.Item(124, 7).Value = 'Dukla Banska Bystrica'
$ws.Cells.Item(124, 8).Value = 2
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 'H'
$ws.Cells.Item(124, 11).Value = 1.666
$ws.Cells.Item(124, 12).Value = 3.75
$ws.Cells.Item(124, 13).Value = 4.5
$ws.Cells.Item(124, 14).Value = 1.615
$ws.Cells.Item(124, 15).Value = 4
$ws.Cells.Item(124, 16).Value = 5.25
$ws.Cells.Item(124, 17).Value = -0.75
$ws.Cells.Item(124, 18).Value = 1.775
$ws.Cells.Item(124, 19).Value = 2.025
$ws.Cells.Item(124, 20).Value = 2.5
$ws.Cells.Item(124, 21).Value = 1.825
$ws.Cells.Item(124, 22).Value = 1.975
$ws.Cells.Item(124, 23).Value = 0.615
$ws.Cells.Item(124, 24).Value = -1
$ws.Cells.Item(124, 25).Value = -1
$ws.Cells.Item(124, 26).Value = 0.7749999999999999
$ws.Cells.Item(124, 27).Value = -1
$ws.Cells.Item(124, 28).Value = -1
$ws.Cells.Item(124, 29).Value = 0.9750000000000001
$ws.Cells.Item(125, 2).Value = 6867488
$ws.Cells.Item(125, 6).Value = 'MSK Zilina'
$ws.Cells.Item(125, 7).Value = 'MFK Zemplin Michalovce'
$ws.Cells.Item(125, 8).Value = 1
$ws.Cells.Item(125, 9).Value = 1
$ws.Cells.Item(125, 10).Value = 'D'
$ws.Cells.Item(125, 11).Value = 1.45
$ws.Cells.Item(125, 12).Value = 4.333
$ws.Cells.Item(125, 13).Value = 5.75
$ws.Cells.Item(125, 14).Value = 1.45
$ws.Cells.Item(125, 15).Value = 4.5
$ws.Cells.Item(125, 16).Value = 6.5
$ws.Cells.Item(125, 17).Value = -1.25
$ws.Cells.Item(125, 18).Value = 2
$ws.Cells.Item(125, 19).Value = 1.8
$ws.Cells.Item(125, 20).Value = 3
$ws.Cells.Item(125, 21).Value = 1.8
$ws.Cells.Item(125, 22).Value = 2
$ws.Cells.Item(125, 23).Value = -1
$ws.Cells.Item(125, 24).Value = 3.5
$ws.Cells.Item(125, 25).Value = -1
$ws.Cells.Item(125, 26).Value = -1
$ws.Cells.Item(125, 27).Value = 0.8
$ws.Cells.Item(125, 28).Value = -1
$ws.Cells.Item(125, 29).Value = 1

# Swap row 128 and row 131 (columns B:AC)
$ws.Cells.Item(128, 2).Value = 6867493
$ws.Cells.Item(128, 6).Value = 'MFK Ruzomberok'
$ws.Cells.Item(128, 7).Value = 'FK Zeleziarne Podbrezova'
$ws.Cells.Item(128, 8).Value = 2
$ws.Cells.Item(128, 9).Value = 1
$ws.Cells.Item(128, 10).Value = 'H'
$ws.Cells.Item(128, 11).Value = 3.3
$ws.Cells.Item(128, 12).Value = 3.3
$ws.Cells.Item(128, 13).Value = 2.2
$ws.Cells.Item(128, 14).Value = 2.9
$ws.Cells.Item(128, 15).Value = 3.4
$ws.Cells.Item(128, 16).Value = 2.4
$ws.Cells.Item(128, 17).Value = 0.25
$ws.Cells.Item(128, 18).Value = 1.75
$ws.Cells.Item(128, 19).Value = 2.05
$ws.Cells.Item(128, 20).Value = 2.5
$ws.Cells.Item(128, 21).Value = 1.9
$ws.Cells.Item(128, 22).Value = 1.9
$ws.Cells.Item(128, 23).Value = 1.9
$ws.Cells.Item(128, 24).Value = -1
$ws.Cells.Item(128, 25).Value = -1
$ws.Cells.Item(128, 26).Value = 0.75
$ws.Cells.Item(128, 27).Value = -1
$ws.Cells.Item(128, 28).Value = 0.8999999999999999
$ws.Cells.Item(128, 29).Value = -1
$ws.Cells.Item(131, 2).Value = 6867491
$ws.Cells.Item(131, 6).Value = 'Dukla Banska Bystrica'
$ws.Cells.Item(131, 7).Value = 'FC Kosice'
$ws.Cells.Item(131, 8).Value = 1
$ws.Cells.Item(131, 9).Value = 1
$ws.Cells.Item(131, 10).Value = 'D'
$ws.Cells.Item(131, 11).Value = 1.571
$ws.Cells.Item(131, 12).Value = 4.2
$ws.Cells.Item(131, 13).Value = 5.25
$ws.Cells.Item(131, 14).Value = 1.571
$ws.Cells.Item(131, 15).Value = 4.2
$ws.Cells.Item(131, 16).Value = 5.5
$ws.Cells.Item(131, 17).Value = -1
$ws.Cells.Item(131, 18).Value = 1.9
$ws.Cells.Item(131, 19).Value = 1.9
$ws.Cells.Item(131, 20).Value = 3
$ws.Cells.Item(131, 21).Value = 2
$ws.Cells.Item(131, 22).Value = 1.8
$ws.Cells.Item(131, 23).Value = -1
$ws.Cells.Item(131, 24).Value = 3.2
$ws.Cells.Item(131, 25).Value = -1
$ws.Cells.Item(131, 26).Value = -1
$ws.Cells.Item(131, 27).Value = 0.8999999999999999
$ws.Cells.Item(131, 28).Value = -1
$ws.Cells.Item(131, 29).Value = 0.8

# Swap row 148 and row 149 (columns B:AC)
$ws.Cells.Item(148, 2).Value = 7911478
$ws.Cells.Item(148, 6).Value = 'FC Vion Zlate Moravce'
$ws.Cells.Item(148, 7).Value = 'FC Kosice'
$ws.Cells.Item(148, 8).Value = 1
$ws.Cells.Item(148, 9).Value = 2
$ws.Cells.Item(148, 10).Value = 'A'
$ws.Cells.Item(148, 11).Value = 2.5
$ws.Cells.Item(148, 12).Value = 3.2
$ws.Cells.Item(148, 13).Value = 2.8
$ws.Cells.Item(148, 14).Value = 2.6
$ws.Cells.Item(148, 15).Value = 3.1
$ws.Cells.Item(148, 16).Value = 2.875
$ws.Cells.Item(148, 17).Value = 0
$ws.Cells.Item(148, 18).Value = 1.775
$ws.Cells.Item(148, 19).Value = 2.025
$ws.Cells.Item(148, 20).Value = 2.25
$ws.Cells.Item(148, 21).Value = 1.8
$ws.Cells.Item(148, 22).Value = 2
$ws.Cells.Item(148, 23).Value = -1
$ws.Cells.Item(148, 24).Value = -1
$ws.Cells.Item(148, 25).Value = 1.875
$ws.Cells.Item(148, 26).Value = -1
$ws.Cells.Item(148, 27).Value = 1.025
$ws.Cells.Item(148, 28).Value = 0.8
$ws.Cells.Item(148, 29).Value = -1
$ws.Cells.Item(149, 2).Value = 7911450
$ws.Cells.Item(149, 6).Value = 'MFK Skalica'
$ws.Cells.Item(149, 7).Value = 'MFK Zemplin Michalovce'
$ws.Cells.Item(149, 8).Value = 0
$ws.Cells.Item(149, 9).Value = 0
$ws.Cells.Item(149, 10).Value = 'D'
$ws.Cells.Item(149, 11).Value = 2.3
$ws.Cells.Item(149, 12).Value = 3.25
$ws.Cells.Item(149, 13).Value = 3.1
$ws.Cells.Item(149, 14).Value = 2.2
$ws.Cells.Item(149, 15).Value = 3.3
$ws.Cells.Item(149, 16).Value = 3.3
$ws.Cells.Item(149, 17).Value = -0.25
$ws.Cells.Item(149, 18).Value = 1.95
$ws.Cells.Item(149, 19).Value = 1.85
$ws.Cells.Item(149, 20).Value = 2.25
$ws.Cells.Item(149, 21).Value = 2
$ws.Cells.Item(149, 22).Value = 1.8
$ws.Cells.Item(149, 23).Value = -1
$ws.Cells.Item(149, 24).Value = 2.3
$ws.Cells.Item(149, 25).Value = -1
$ws.Cells.Item(149, 26).Value = -0.5
$ws.Cells.Item(149, 27).Value = 0.425
$ws.Cells.Item(149, 28).Value = -1
$ws.Cells.Item(149, 29).Value = 0.8

# Update odds values in rows 152-157
# Row 152
$ws.Cells.Item(152, 14).Value = 5.25
$ws.Cells.Item(152, 15).Value = 4
$ws.Cells.Item(152, 16).Value = 1.65
$ws.Cells.Item(152, 18).Value = 1.8
$ws.Cells.Item(152, 19).Value = 2
$ws.Cells.Item(152, 21).Value = 1.975
$ws.Cells.Item(152, 22).Value = 1.825

# Row 153
$ws.Cells.Item(153, 14).Value = 3.8
$ws.Cells.Item(153, 15).Value = 3.3
$ws.Cells.Item(153, 16).Value = 2.05
$ws.Cells.Item(153, 18).Value = 1.775
$ws.Cells.Item(153, 19).Value = 2.025

# Row 154
$ws.Cells.Item(154, 14).Value = 3
$ws.Cells.Item(154, 15).Value = 3.3
$ws.Cells.Item(154, 16).Value = 2.375
$ws.Cells.Item(154, 18).Value = 1.8
$ws.Cells.Item(154, 19).Value = 2
$ws.Cells.Item(154, 21).Value = 2
$ws.Cells.Item(154, 22).Value = 1.8

# Row 155
$ws.Cells.Item(155, 14).Value = 3.4
$ws.Cells.Item(155, 16).Value = 2.15
$ws.Cells.Item(155, 18).Value = 2
$ws.Cells.Item(155, 19).Value = 1.8
$ws.Cells.Item(155, 21).Value = 1.8
$ws.Cells.Item(155, 22).Value = 2

# Row 156
$ws.Cells.Item(156, 14).Value = 1.8
$ws.Cells.Item(156, 15).Value = 3.5
$ws.Cells.Item(156, 16).Value = 4.75
$ws.Cells.Item(156, 20).Value = 2.5
$ws.Cells.Item(156, 21).Value = 2.025
$ws.Cells.Item(156, 22).Value = 1.775

# Row 157
$ws.Cells.Item(157, 16).Value = 4.5
